# "Report formulas & format"
# - Disable iterative calculation (workbook calc options).
# - Clear the placeholder "." value out of A3 (template cell left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off iterative calculation (File > Options > Formulas > Enable iterative
# calculation), restoring the engine defaults for iteration count / max change.
$excel.Iteration = $false
$excel.MaxIterations = 100
$excel.MaxChange = 0.001

# Clear the leftover "." placeholder text from the template's third row.
$ws.Range("A3").Value = $null
